$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "components" sentence: the Zumo "shield" became a Zumo "Robot",
#    and a stray double-space before "a X-Bee Shield" was cleaned up.
# ------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute(
    "a Zumo shield, a Microsoft Kinect, an Arduino Uno,  a X-Bee Shield",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "a Zumo Robot, a Microsoft Kinect, an Arduino Uno, a X-Bee Shield",
    2)
Write-Host "Components sentence updated:" $found

# ------------------------------------------------------------------
# 2) The "_GoBack" bookmark moves from the end of the document
#    (it used to sit right after the stray "4." paragraph near the
#    end) to sit right after the word "Robot" in the sentence we
#    just edited, i.e. right before ", a Microsoft Kinect...".
# ------------------------------------------------------------------
$existing = $d.Bookmarks.Item("_GoBack")
if ($existing -ne $null) {
    $existing.Delete()
}

$r2 = $d.Content
$found2 = $r2.Find.Execute(
    "Robot, a Microsoft Kinect",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "", 0)
Write-Host "Located new bookmark anchor:" $found2

$insertPoint = $r2.Start + 5   # length of "Robot" -> right after it, before the comma
$bmRange = $d.Range($insertPoint, $insertPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
Write-Host "Bookmark moved"
